$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "path to pic" column (H2:H7) stores Windows-style backslash paths
# ("PicturesOrig\Xxx.jpg"); change the path separator to a forward slash
# ("PicturesOrig/Xxx.jpg") for every row.
$ws.Range("H2:H7").Replace("\", "/")

# Move/restore the saved selection to H7 (within the sheet's used range,
# A1:H7) instead of the stale H13 reference.
$ws.Range("H7").Select()
